$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

$ws_ALC.Range("H41").Value = 1091.6
$ws_ALC.Range("J41").Value = 1239.5
$ws_ALC.Range("L41").Value = 1239.5
$ws_ALC.Range("N41").Value = -2119.5

$ws_ALC.Range("H55").Value = 130.66667
$ws_ALC.Range("J55").Value = 78.28570999999999
$ws_ALC.Range("L55").Value = 78.28570999999999
$ws_ALC.Range("N55").Value = -506.28571

$ws_ALC.Range("H103").Value = 5409.5
$ws_ALC.Range("I103").Value = 0
$ws_ALC.Range("K103").Value = 0
$ws_ALC.Range("M103").ClearContents()

$ws_ALC.Range("H113").Value = 10274.6875
$ws_ALC.Range("J113").Value = 6235.909
$ws_ALC.Range("L113").Value = 6235.909
$ws_ALC.Range("N113").Value = -12743.909

$ws_ALC.Range("H132").Value = 10790.315
$ws_ALC.Range("I132").Value = 3078.6924
$ws_ALC.Range("K132").Value = 9236.0772
$ws_ALC.Range("M132").Value = -6706.0772

$ws_ALC.Range("H135").Value = 5686.7
$ws_ALC.Range("I135").Value = 6124.1113
$ws_ALC.Range("K135").Value = 55117.00169999999
$ws_ALC.Range("M135").Value = -52582.00169999999

$ws_ARM.Range("H61").Value = 605932.6
$ws_ARM.Range("I61").Value = 3712.2307
$ws_ARM.Range("K61").Value = 3712.2307
$ws_ARM.Range("M61").Value = -3500.2307

$ws_ARM.Range("H122").Value = 834865.25
$ws_ARM.Range("I122").Value = 1001468.8
$ws_ARM.Range("K122").Value = 3004406.4
$ws_ARM.Range("M122").Value = -3001956.4

$ws_ARM.Range("H132").Value = 811780.5600000001
$ws_ARM.Range("I132").Value = 5120.7188
$ws_ARM.Range("K132").Value = 15362.1564
$ws_ARM.Range("M132").Value = -12832.1564

$ws_ARM.Range("H136").Value = 605932.6
$ws_ARM.Range("I136").Value = 3712.2307
$ws_ARM.Range("K136").Value = 11136.6921
$ws_ARM.Range("M136").Value = -8586.6921

$ws_BSM.Range("H86").Value = 2000
$ws_BSM.Range("I86").Value = 2000
$ws_BSM.Range("K86").Value = 2000
$ws_BSM.Range("M86").Value = -877

$ws_BSM.Range("H89").Value = 2000
$ws_BSM.Range("I89").Value = 2000
$ws_BSM.Range("K89").Value = 10000
$ws_BSM.Range("M89").Value = -4384

$ws_BSM.Range("H94").Value = 734.7857
$ws_BSM.Range("I94").Value = 734.0454999999999
$ws_BSM.Range("K94").Value = 734.0454999999999
$ws_BSM.Range("M94").Value = -283.0454999999999

$ws_BSM.Range("H107").Value = 7752.5
$ws_BSM.Range("I107").Value = 8670
$ws_BSM.Range("J107").Value = 5000
$ws_BSM.Range("K107").Value = 8670
$ws_BSM.Range("L107").Value = 5000
$ws_BSM.Range("M107").Value = -6750
$ws_BSM.Range("N107").Value = -8840

$ws_CRP.Range("H7").Value = 832.7
$ws_CRP.Range("I7").Value = 71.59999999999999
$ws_CRP.Range("J7").Value = 1593.8
$ws_CRP.Range("K7").Value = 71.59999999999999
$ws_CRP.Range("L7").Value = 1593.8
$ws_CRP.Range("M7").Value = 41.40000000000001
$ws_CRP.Range("N7").Value = -1819.8

$ws_CRP.Range("H31").Value = 90936.484
$ws_CRP.Range("J31").Value = 43222.7
$ws_CRP.Range("L31").Value = 43222.7
$ws_CRP.Range("N31").Value = -43812.7

$ws_CRP.Range("H34").Value = 90936.484
$ws_CRP.Range("J34").Value = 43222.7
$ws_CRP.Range("L34").Value = 43222.7
$ws_CRP.Range("N34").Value = -43626.7

$ws_CRP.Range("H58").Value = 13014.172
$ws_CRP.Range("I58").Value = 4909.25
$ws_CRP.Range("K58").Value = 4909.25
$ws_CRP.Range("M58").Value = -4706.25

$ws_CRP.Range("H132").Value = 2959876
$ws_CRP.Range("I132").Value = 4496.8335
$ws_CRP.Range("K132").Value = 13490.5005
$ws_CRP.Range("M132").Value = -10960.5005

$ws_CRP.Range("H134").Value = 6444.2144
$ws_CRP.Range("I134").Value = 2757.0435
$ws_CRP.Range("K134").Value = 8271.130500000001
$ws_CRP.Range("M134").Value = -5736.130500000001

$ws_CRP.Range("H136").Value = 13014.172
$ws_CRP.Range("I136").Value = 4909.25
$ws_CRP.Range("K136").Value = 14727.75
$ws_CRP.Range("M136").Value = -12177.75

$ws_CUL.Range("H37").Value = 500044000
$ws_CUL.Range("J37").Value = 500044000
$ws_CUL.Range("L37").Value = 1500132000
$ws_CUL.Range("N37").Value = -1500132224

$ws_CUL.Range("H132").Value = 1895445.5
$ws_CUL.Range("I132").Value = 1610
$ws_CUL.Range("J132").Value = 5051838
$ws_CUL.Range("K132").Value = 14490
$ws_CUL.Range("L132").Value = 45466542
$ws_CUL.Range("M132").Value = -11960
$ws_CUL.Range("N132").Value = -45471602

$ws_CUL.Range("H134").Value = 5037.795
$ws_CUL.Range("I134").Value = 1851.4117
$ws_CUL.Range("K134").Value = 5554.2351
$ws_CUL.Range("M134").Value = -484.2350999999999

$ws_GSM.Range("H2").Value = 153
$ws_GSM.Range("I2").Value = 251.66667
$ws_GSM.Range("K2").Value = 251.66667
$ws_GSM.Range("M2").Value = -138.66667

$ws_GSM.Range("H41").Value = 4233.6665
$ws_GSM.Range("I41").Value = 5356.75
$ws_GSM.Range("K41").Value = 5356.75
$ws_GSM.Range("M41").Value = -5001.75

$ws_GSM.Range("H80").Value = 2600.3333
$ws_GSM.Range("I80").Value = 2595
$ws_GSM.Range("K80").Value = 2595
$ws_GSM.Range("M80").Value = -1597

$ws_GSM.Range("H83").Value = 2600.3333
$ws_GSM.Range("I83").Value = 2595
$ws_GSM.Range("K83").Value = 12975
$ws_GSM.Range("M83").Value = -7983

$ws_GSM.Range("H102").Value = 6915.8335
$ws_GSM.Range("I102").Value = 7299
$ws_GSM.Range("K102").Value = 7299
$ws_GSM.Range("M102").Value = -5677

$ws_GSM.Range("H107").Value = 1271.8572
$ws_GSM.Range("J107").Value = 400
$ws_GSM.Range("L107").Value = 400
$ws_GSM.Range("N107").Value = -4240

$ws_GSM.Range("H113").Value = 4816.476
$ws_GSM.Range("I113").Value = 4977.4
$ws_GSM.Range("J113").Value = 4414.1665
$ws_GSM.Range("K113").Value = 4977.4
$ws_GSM.Range("L113").Value = 4414.1665
$ws_GSM.Range("M113").Value = -2807.4
$ws_GSM.Range("N113").Value = -8754.166499999999

$ws_GSM.Range("H122").Value = 1337.5333
$ws_GSM.Range("I122").Value = 968.1818
$ws_GSM.Range("K122").Value = 2904.5454
$ws_GSM.Range("M122").Value = -454.5454

$ws_GSM.Range("H132").Value = 659366.4
$ws_GSM.Range("I132").Value = 4068.2354
$ws_GSM.Range("J132").Value = 1897151.8
$ws_GSM.Range("K132").Value = 12204.7062
$ws_GSM.Range("L132").Value = 5691455.4
$ws_GSM.Range("M132").Value = -9674.706200000001
$ws_GSM.Range("N132").Value = -5696515.4

$ws_LTW.Range("H7").Value = 6253.3657
$ws_LTW.Range("I7").Value = 6861.08
$ws_LTW.Range("K7").Value = 6861.08
$ws_LTW.Range("M7").Value = -6749.08

$ws_LTW.Range("H16").Value = 1301.7858
$ws_LTW.Range("I16").Value = 1211.931
$ws_LTW.Range("J16").Value = 1502.2307
$ws_LTW.Range("K16").Value = 1211.931
$ws_LTW.Range("L16").Value = 1502.2307
$ws_LTW.Range("M16").Value = -1041.931
$ws_LTW.Range("N16").Value = -1842.2307

$ws_LTW.Range("H22").Value = 933.8
$ws_LTW.Range("I22").Value = 1299.5
$ws_LTW.Range("K22").Value = 1299.5
$ws_LTW.Range("M22").Value = -1004.5

$ws_LTW.Range("H27").Value = 933.8
$ws_LTW.Range("I27").Value = 1299.5
$ws_LTW.Range("K27").Value = 1299.5
$ws_LTW.Range("M27").Value = -1192.5

$ws_LTW.Range("H46").Value = 3030.1428
$ws_LTW.Range("I46").Value = 837
$ws_LTW.Range("J46").Value = 3546.1765
$ws_LTW.Range("K46").Value = 837
$ws_LTW.Range("L46").Value = 3546.1765
$ws_LTW.Range("M46").Value = -649
$ws_LTW.Range("N46").Value = -3922.1765

$ws_LTW.Range("H55").Value = 597.5833
$ws_LTW.Range("I55").Value = 464.83334
$ws_LTW.Range("K55").Value = 464.83334
$ws_LTW.Range("M55").Value = -291.83334

$ws_LTW.Range("H61").Value = 2798.05
$ws_LTW.Range("I61").Value = 2379.7646
$ws_LTW.Range("J61").Value = 5168.3335
$ws_LTW.Range("K61").Value = 2379.7646
$ws_LTW.Range("L61").Value = 5168.3335
$ws_LTW.Range("M61").Value = -2177.7646
$ws_LTW.Range("N61").Value = -5572.3335

$ws_LTW.Range("H100").Value = 2976.3914
$ws_LTW.Range("J100").Value = 3540.5715
$ws_LTW.Range("L100").Value = 3540.5715
$ws_LTW.Range("N100").Value = -4622.5715

$ws_LTW.Range("H113").Value = 2798.05
$ws_LTW.Range("I113").Value = 2379.7646
$ws_LTW.Range("J113").Value = 5168.3335
$ws_LTW.Range("K113").Value = 2379.7646
$ws_LTW.Range("L113").Value = 5168.3335
$ws_LTW.Range("M113").Value = -209.7646
$ws_LTW.Range("N113").Value = -9508.333500000001

$ws_LTW.Range("H126").Value = 6253.3657
$ws_LTW.Range("I126").Value = 6861.08
$ws_LTW.Range("K126").Value = 20583.24
$ws_LTW.Range("M126").Value = -18113.24

$ws_LTW.Range("H132").Value = 1000779.06
$ws_LTW.Range("I132").Value = 4242.5264
$ws_LTW.Range("K132").Value = 12727.5792
$ws_LTW.Range("M132").Value = -10197.5792

$ws_WVR.Range("H107").Value = 724.8461
$ws_WVR.Range("I107").Value = 743.125
$ws_WVR.Range("J107").Value = 695.6
$ws_WVR.Range("K107").Value = 2229.375
$ws_WVR.Range("L107").Value = 2086.8
$ws_WVR.Range("M107").Value = -309.375
$ws_WVR.Range("N107").Value = -5926.8

$ws_WVR.Range("H132").Value = 426612.72
$ws_WVR.Range("J132").Value = 1000756.8
$ws_WVR.Range("L132").Value = 3002270.4
$ws_WVR.Range("N132").Value = -3007330.4

$ws_WVR.Range("H136").Value = 464017.8
$ws_WVR.Range("I136").Value = 2658.8
$ws_WVR.Range("J136").Value = 818909.4
$ws_WVR.Range("K136").Value = 7976.400000000001
$ws_WVR.Range("L136").Value = 2456728.2
$ws_WVR.Range("M136").Value = -5426.400000000001
$ws_WVR.Range("N136").Value = -2461828.2

